# Assignment07_Writeup.docx edit
#
# 1) "Github:" line -> split the trailing run into a leading-space run and a
#    URL run whose repo name gains a "-Module07" suffix.
# 2) Two standalone "Go" (T-SQL batch separator) paragraphs that are each
#    stored as two runs ("G" + "o") get collapsed into a single "Go" run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a Range's content with a single run carrying the given
# run-properties fragment (rPrXml, e.g. "<w:rPr>...</w:rPr>" or "") and text.
# Using Range.InsertXML (instead of Range.Text = ...) forces Word to keep the
# inserted text in its own run rather than silently re-merging it into a
# neighbouring run that happens to share the same formatting.
# ---------------------------------------------------------------------------
function Set-RangeAsSingleRun($range, [string]$rPrXml, [string]$text) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    if ($text.StartsWith(" ") -or $text.EndsWith(" ")) {
        $tTag = '<w:t xml:space="preserve">' + $escaped + '</w:t>'
    } else {
        $tTag = '<w:t>' + $escaped + '</w:t>'
    }
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $rPrXml + $tTag + '</w:r></w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) Github URL: " https://github.com/guillermo-dominguez/DBFoundations"
#    -> " " (run, keeps its original run identity/attributes) +
#       "https://github.com/guillermo-dominguez/DBFoundations-Module07" (new run)
# ---------------------------------------------------------------------------
$oldUrl = "https://github.com/guillermo-dominguez/DBFoundations"
$newUrl = "https://github.com/guillermo-dominguez/DBFoundations-Module07"
$urlRPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$full = $d.Content.Text
$pos = $full.IndexOf(" " + $oldUrl)
if ($pos -ge 0) {
    $urlStart = $pos + 1
    $urlEnd = $urlStart + $oldUrl.Length

    # First, nudge the font of the URL-only sub-range so Word is forced to
    # split it off into its own run, leaving the untouched leading-space run
    # (and its rsid/other <w:r> attributes) exactly as it was.
    $splitRange = $d.Range($urlStart, $urlEnd)
    $splitRange.Font.Name = "Courier New"

    # Now overwrite that freshly-split run with the final text/formatting.
    $urlRange = $d.Range($urlStart, $urlEnd)
    Set-RangeAsSingleRun $urlRange $urlRPr $newUrl
}

# ---------------------------------------------------------------------------
# 2) Merge the two "G" + "o" runs (Consolas, blue, sz 19) into a single "Go"
#    run, wherever that exact two-run split occurs.
# ---------------------------------------------------------------------------
$goRPr = '<w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr>'

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $pRange = $p.Range
    $pStart = $pRange.Start
    $pEnd = $pRange.End
    $bodyRange = $d.Range($pStart, $pEnd - 1)
    $text = $bodyRange.Text
    if ($text.Equals("Go")) {
        Set-RangeAsSingleRun $bodyRange $goRPr "Go"
    }
}
